# Refresh the cryptos list with the latest scraped price/volume snapshot.
# Mirrors the upstream GitHub Actions job that re-runs the scraper and
# commits whatever cells changed since the previous run (some rows only
# change in the Volume(1h) column, row 26/27 swap rank because LidoDAOToken
# overtook EthereumClassic).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) cells hold plain text in the source sheet (e.g. the
# thousands-grouped "27.113.78") so values that *look* numeric need the
# Text number format forced before the write - otherwise COM Value
# assignment auto-converts them to real numbers. ClearFormats() afterwards
# drops the cell back to the workbook default style (so no stray format
# delta is left behind) while leaving the stored value as text.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.113.78"
$ws.Range("E2").Value = "  -0.97%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.780.50"
$ws.Range("E3").Value = "  -1.84%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "336.53"
$ws.Range("E5").Value = "  -2.44%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.10%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.3821"
$ws.Range("E7").Value = "  +0.33%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.3409"
$ws.Range("E8").Value = "  -2.30%  "

# Row 9 - OKB
Set-TextValue $ws.Range("D9") "47.96"
$ws.Range("E9").Value = "  -1.84%  "

# Row 10 - Polygon
Set-TextValue $ws.Range("D10") "1.185"

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.07444"
$ws.Range("E11").Value = "  -3.65%  "

# Row 12 - BinanceUSD
Set-TextValue $ws.Range("D12") "1.002"
$ws.Range("E12").Value = "  -0.10%  "

# Row 13 - Solana
Set-TextValue $ws.Range("D13") "21.60"
$ws.Range("E13").Value = "  -2.15%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.419"
$ws.Range("E14").Value = "  -2.97%  "

# Row 15 - WrappedEther
$ws.Range("D15").Value = "1.781.38"
$ws.Range("E15").Value = "  -1.87%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "7.104"
$ws.Range("E16").Value = "  -1.62%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.00001090"
$ws.Range("E17").Value = "  -2.36%  "

# Row 18 - TRON
Set-TextValue $ws.Range("D18") "0.06648"
$ws.Range("E18").Value = "  -1.26%  "

# Row 19 - Litecoin
Set-TextValue $ws.Range("D19") "83.35"
$ws.Range("E19").Value = "  -3.02%  "

# Row 20 - Dai
Set-TextValue $ws.Range("D20") "1.002"
$ws.Range("E20").Value = "  +0.03%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.511"
$ws.Range("E21").Value = "  -1.10%  "

# Row 22 - Avalanche
Set-TextValue $ws.Range("D22") "17.35"
$ws.Range("E22").Value = "  -1.53%  "

# Row 23 - WrappedBTC
$ws.Range("D23").Value = "27.103.41"
$ws.Range("E23").Value = "  -1.05%  "

# Row 24 - Cosmos
$ws.Range("E24").Value = "  -7.71%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.373"
$ws.Range("E25").Value = "  -3.95%  "

# Row 26 - LidoDAOToken (was EthereumClassic - rows 26/27 swapped)
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D26") "2.494"
$ws.Range("E26").Value = "  -6.46%  "

# Row 27 - EthereumClassic (was LidoDAOToken - rows 26/27 swapped)
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D27") "21.11"
$ws.Range("E27").Value = "  -3.91%  "

# Row 28 - ImmutableX
Set-TextValue $ws.Range("D28") "1.447"
$ws.Range("E28").Value = "  -1.30%  "

# Row 29 - Monero
Set-TextValue $ws.Range("D29") "154.04"
$ws.Range("E29").Value = "  -0.29%  "

# Row 30 - WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.984.66"
$ws.Range("E30").Value = "  -1.72%  "

# Row 31 - BitcoinCash
Set-TextValue $ws.Range("D31") "133.82"
$ws.Range("E31").Value = "  -1.43%  "

# Row 32 - HuobiToken
Set-TextValue $ws.Range("D32") "3.972"
$ws.Range("E32").Value = "  -1.41%  "

# Row 33 - Filecoin
Set-TextValue $ws.Range("D33") "5.997"
$ws.Range("E33").Value = "  -4.95%  "

# Row 34 - Stellar
Set-TextValue $ws.Range("D34") "0.08656"
$ws.Range("E34").Value = "  -1.29%  "

# Row 35 - Aptos
Set-TextValue $ws.Range("D35") "12.99"
$ws.Range("E35").Value = "  -6.57%  "

# Row 36 - WEMIXTOKEN
$ws.Range("E36").Value = "  -4.06%  "

# Row 37 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D37") "5.379"
$ws.Range("E37").Value = "  -4.19%  "

# Row 38 - TheSandbox
Set-TextValue $ws.Range("D38") "0.6811"
$ws.Range("E38").Value = "  -2.33%  "

# Row 39 - Hedera
Set-TextValue $ws.Range("D39") "0.06318"
$ws.Range("E39").Value = "  -2.38%  "

# Row 40 - VeChain
Set-TextValue $ws.Range("D40") "0.02326"
$ws.Range("E40").Value = "  -3.25%  "

# Row 41 - Algorand
Set-TextValue $ws.Range("D41") "0.2172"
$ws.Range("E41").Value = "  -4.37%  "

# Row 42 - TrustWalletToken
Set-TextValue $ws.Range("D42") "1.237"
$ws.Range("E42").Value = "  -5.26%  "

# Row 43 - FraxShare
Set-TextValue $ws.Range("D43") "8.428"
$ws.Range("E43").Value = "  -5.64%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "14.26"
$ws.Range("E44").Value = "  -2.96%  "

# Row 45 - Frax
$ws.Range("E45").Value = "  +0.02%  "

# Row 46 - Decentraland
Set-TextValue $ws.Range("D46") "0.6395"
$ws.Range("E46").Value = "  -1.96%  "

# Row 47 - PancakeSwap
$ws.Range("E47").Value = "  -4.63%  "

# Row 48 - NEARProtocol
Set-TextValue $ws.Range("D48") "2.150"
$ws.Range("E48").Value = "  -1.26%  "

# Row 49 - Quant
Set-TextValue $ws.Range("D49") "130.95"
$ws.Range("E49").Value = "  -1.22%  "

# Row 50 - Cronos
Set-TextValue $ws.Range("D50") "0.07083"
$ws.Range("E50").Value = "  -3.24%  "

# Row 51 - Aave
Set-TextValue $ws.Range("D51") "78.46"
$ws.Range("E51").Value = "  -2.44%  "
